$wb = $excel.ActiveWorkbook
$wsShort = $wb.Worksheets.Item("Short Term")
$wsMedium = $wb.Worksheets.Item("Medium Term")

# ---- Sheet "Short Term" updates (rows 110-117) ----
$wsShort.Cells.Item(110, 2).Value = 8.31
$wsShort.Cells.Item(110, 3).Value = 14.95
$wsShort.Cells.Item(110, 4).Value = -14.15
$wsShort.Cells.Item(110, 6).Value = 22.01

$wsShort.Cells.Item(111, 2).Value = -3.96
$wsShort.Cells.Item(111, 3).Value = -2.37
$wsShort.Cells.Item(111, 4).Value = 7.86

$wsShort.Cells.Item(112, 2).Value = 5.22
$wsShort.Cells.Item(112, 3).Value = 2.75
$wsShort.Cells.Item(112, 4).Value = 5.59
$wsShort.Cells.Item(112, 5).Value = 0.05
$wsShort.Cells.Item(112, 6).Value = 2.31
$wsShort.Cells.Item(112, 7).Value = -0.63

$wsShort.Cells.Item(113, 2).Value = -0.48
$wsShort.Cells.Item(113, 3).Value = 6.03
$wsShort.Cells.Item(113, 4).Value = -4.75

$wsShort.Cells.Item(114, 2).Value = -2.03
$wsShort.Cells.Item(114, 3).Value = -10.57
$wsShort.Cells.Item(114, 4).Value = 0.24

$wsShort.Cells.Item(115, 2).Value = -2.78
$wsShort.Cells.Item(115, 3).Value = -0.87
$wsShort.Cells.Item(115, 4).Value = 6.54

$wsShort.Cells.Item(116, 2).Value = 4.83
$wsShort.Cells.Item(116, 3).Value = 3.25
$wsShort.Cells.Item(116, 4).Value = 3.41

$wsShort.Cells.Item(117, 2).Value = 0.82
$wsShort.Cells.Item(117, 4).Value = 6.14
$wsShort.Cells.Item(117, 5).Value = 8.39
$wsShort.Cells.Item(117, 6).Value = 7.44
$wsShort.Cells.Item(117, 7).Value = 11.11

# New row 118 on "Short Term" (A118 gets the same date style as A117; B118/E118 stay empty)
$wsShort.Cells.Item(117, 1).Copy($wsShort.Cells.Item(118, 1))
$wsShort.Cells.Item(118, 1).Value = 45536
$wsShort.Cells.Item(2, 2).Copy($wsShort.Cells.Item(118, 2))
$wsShort.Cells.Item(118, 3).Value = 22.64
$wsShort.Cells.Item(118, 4).Value = -12.11
$wsShort.Cells.Item(2, 2).Copy($wsShort.Cells.Item(118, 5))
$wsShort.Cells.Item(118, 6).Value = 40.35
$wsShort.Cells.Item(118, 7).Value = -5.42

# ---- Sheet "Medium Term" updates (rows 96-103) ----
$wsMedium.Cells.Item(96, 2).Value = 1.71
$wsMedium.Cells.Item(96, 3).Value = -5.84
$wsMedium.Cells.Item(96, 4).Value = -4.21

$wsMedium.Cells.Item(97, 2).Value = 2.32
$wsMedium.Cells.Item(97, 3).Value = -3.53
$wsMedium.Cells.Item(97, 4).Value = -4.97

$wsMedium.Cells.Item(98, 2).Value = 6.02
$wsMedium.Cells.Item(98, 3).Value = 0.43
$wsMedium.Cells.Item(98, 4).Value = -4.06

$wsMedium.Cells.Item(99, 2).Value = 8.03
$wsMedium.Cells.Item(99, 3).Value = 4.95
$wsMedium.Cells.Item(99, 4).Value = -0.95

$wsMedium.Cells.Item(100, 2).Value = 12.39
$wsMedium.Cells.Item(100, 3).Value = 7.51
$wsMedium.Cells.Item(100, 4).Value = 0.38

$wsMedium.Cells.Item(101, 3).Value = 7.21
$wsMedium.Cells.Item(101, 4).Value = -0.73

$wsMedium.Cells.Item(102, 3).Value = 6.09
$wsMedium.Cells.Item(102, 4).Value = 0.03

$wsMedium.Cells.Item(103, 2).Value = 3.74
$wsMedium.Cells.Item(103, 3).Value = 8.08
$wsMedium.Cells.Item(103, 4).Value = 2.25

# New row 104 on "Medium Term" (A104 gets the same date style as A103)
$wsMedium.Cells.Item(103, 1).Copy($wsMedium.Cells.Item(104, 1))
$wsMedium.Cells.Item(104, 1).Value = 45536
$wsMedium.Cells.Item(104, 2).Value = 20.74
$wsMedium.Cells.Item(104, 3).Value = 14.42
$wsMedium.Cells.Item(104, 4).Value = 7.31
